# fpt-report1-template.docx edit script
# - update cover date text
# - add a page-number footer to the (single) section
# - retune paragraph/heading/list/table-style typography to Calibri

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Cover page date: "Hanoi, August 2019" -> "Hanoi, January 2026"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Hanoi, August 2019", $true, $false, $false, $false, `
                         $false, $true, 1, $false, "Hanoi, January 2026", 2) | Out-Null

# ---------------------------------------------------------------------------
# 2. Footer with centered PAGE field, linked to the only section
# ---------------------------------------------------------------------------
$footerXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/temp.xml" pkg:contentType="text/xml">
    <pkg:xmlData>
      <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:pPr>
          <w:jc w:val="center"/>
        </w:pPr>
        <w:r>
          <w:fldChar w:fldCharType="begin"/>
        </w:r>
        <w:r>
          <w:instrText>PAGE</w:instrText>
        </w:r>
        <w:r>
          <w:fldChar w:fldCharType="end"/>
        </w:r>
      </w:p>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$footer = $d.Sections.Item(1).Footers.Item(1)
$footer.Range.InsertXML($footerXml) | Out-Null
$footerPara = $footer.Range.Paragraphs.Item(1)
$footerPara.Style = "Footer"
$footerPara.Alignment = 1

# ---------------------------------------------------------------------------
# 3. Typography pass across the template styles
# ---------------------------------------------------------------------------

# Normal: base font + default paragraph spacing
$normal = $d.Styles.Item("Normal")
$normal.Font.Name = "Calibri"
$normal.Font.NameFarEast = "Calibri"
$normal.ParagraphFormat.SpaceBefore = 0
$normal.ParagraphFormat.SpaceAfter = 6

# Heading 1: 18pt before / 6pt after
$h1 = $d.Styles.Item("Heading1")
$h1.Font.Name = "Calibri"
$h1.Font.NameFarEast = "Calibri"
$h1.ParagraphFormat.SpaceBefore = 18
$h1.ParagraphFormat.SpaceAfter = 6

# Heading 2: 12pt before / 6pt after
$h2 = $d.Styles.Item("Heading2")
$h2.Font.Name = "Calibri"
$h2.Font.NameFarEast = "Calibri"
$h2.ParagraphFormat.SpaceBefore = 12
$h2.ParagraphFormat.SpaceAfter = 6

# Heading 3: 10pt before / 4pt after
$h3 = $d.Styles.Item("Heading3")
$h3.Font.Name = "Calibri"
$h3.Font.NameFarEast = "Calibri"
$h3.ParagraphFormat.SpaceBefore = 10
$h3.ParagraphFormat.SpaceAfter = 4

# List Paragraph: hanging-friendly spacing + Calibri
$listPara = $d.Styles.Item("ListParagraph")
$listPara.Font.Name = "Calibri"
$listPara.Font.NameFarEast = "Calibri"
$listPara.ParagraphFormat.SpaceBefore = 0
$listPara.ParagraphFormat.SpaceAfter = 3

# Table Text small: 11pt body text in tables, keep cs=Times New Roman
$tableText = $d.Styles.Item("TableTextsmall")
$tableText.Font.Name = "Calibri"
$tableText.Font.NameFarEast = "Calibri"

# Table Head: bold table-header labels, keep cs=Times New Roman
$tableHead = $d.Styles.Item("TableHead")
$tableHead.Font.Name = "Calibri"
$tableHead.Font.NameFarEast = "Calibri"

Write-Output "edit applied"
